# Refresh the "cryptos" price list with the latest Price / Volume(1h) figures.
# Notes:
#  - D/E columns hold plain text (prices such as "26.883.34" are not valid
#    numbers and must stay text; percentages keep their padding spaces).
#  - A handful of D-column values look like ordinary decimals (e.g. "205.57").
#    Excel would otherwise auto-convert those to numbers, so they are entered
#    with a leading apostrophe (forces text) and then the cell style is reset
#    to "Normal" so no stray number-format/style gets attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.883.34"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.544.40"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'205.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'21.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.765.13"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.543.67"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'0.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "26.861.74"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'61.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'213.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "0.0₃0681"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'7.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "'152.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").Value = "'14.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "1.354.92"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'1.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'0.965"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.03%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "'0.803"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'5.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").Value = "'0.989"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'63.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "1.679.37"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "'85.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "'0.0946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "
